# Update the academy admission year from 2022 to 2023.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2022", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2023", 2)
